$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value2 = '44.091.68'
$ws.Cells.Item(2, 5).Value2 = '  +5.35%  '

# Row 3
$ws.Cells.Item(3, 4).Value2 = '2.260.99'
$ws.Cells.Item(3, 5).Value2 = '  +2.59%  '

# Row 4
$ws.Cells.Item(4, 5).Value2 = '  -0.10%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = '230.34'
$ws.Cells.Item(5, 5).Value2 = '  +0.07%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value2 = '0.633'
$ws.Cells.Item(6, 5).Value2 = '  +3.36%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value2 = '63.70'
$ws.Cells.Item(7, 5).Value2 = '  +5.38%  '

# Row 8
$ws.Cells.Item(8, 5).Value2 = '  -0.05%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value2 = '0.450'
$ws.Cells.Item(9, 5).Value2 = '  +12.54%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value2 = '0.104'
$ws.Cells.Item(10, 5).Value2 = '  +16.07%  '

# Row 11
$ws.Cells.Item(11, 5).Value2 = '  -0.47%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value2 = '26.33'
$ws.Cells.Item(12, 5).Value2 = '  +19.89%  '

# Row 13
$ws.Cells.Item(13, 5).Value2 = '  +2.26%  '

# Row 14
$ws.Cells.Item(14, 4).Value2 = '2.594.56'
$ws.Cells.Item(14, 5).Value2 = '  +2.47%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value2 = '15.69'
$ws.Cells.Item(15, 5).Value2 = '  +2.24%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value2 = '6.08'
$ws.Cells.Item(16, 5).Value2 = '  +9.64%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value2 = '0.835'
$ws.Cells.Item(17, 5).Value2 = '  +5.40%  '

# Row 18
$ws.Cells.Item(18, 4).Value2 = '2.254.57'
$ws.Cells.Item(18, 5).Value2 = '  +2.56%  '

# Row 19
$ws.Cells.Item(19, 4).Value2 = '43.943.97'
$ws.Cells.Item(19, 5).Value2 = '  +4.97%  '

# Row 20
$ws.Cells.Item(20, 5).Value2 = '  +8.82%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value2 = '73.44'
$ws.Cells.Item(21, 5).Value2 = '  +2.29%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value2 = '6.03'
$ws.Cells.Item(22, 5).Value2 = '  -2.13%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value2 = '251.47'
$ws.Cells.Item(23, 5).Value2 = '  +4.00%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value2 = '1.00'
$ws.Cells.Item(24, 5).Value2 = '  +0.10%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = '2.43'
$ws.Cells.Item(25, 5).Value2 = '  +1.89%  '

# Row 26
$ws.Cells.Item(26, 5).Value2 = '  -1.67%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value2 = '3.36'
$ws.Cells.Item(27, 5).Value2 = '  +26.79%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value2 = '10.08'
$ws.Cells.Item(28, 5).Value2 = '  +5.37%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value2 = '172.11'
$ws.Cells.Item(29, 5).Value2 = '  +1.96%  '

# Row 30: 'Kaspa' -> 'EthereumClassic'
$ws.Cells.Item(30, 2).Value2 = 'EthereumClassic'
$ws.Cells.Item(30, 3).Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value2 = '20.78'
$ws.Cells.Item(30, 5).Value2 = '  +3.38%  '

# Row 31: 'EthereumClassic' -> 'Kaspa'
$ws.Cells.Item(31, 2).Value2 = 'Kaspa'
$ws.Cells.Item(31, 3).Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value2 = '0.137'
$ws.Cells.Item(31, 5).Value2 = '  -1.74%  '

# Row 32
$ws.Cells.Item(32, 5).Value2 = '  -2.55%  '

# Row 33
$ws.Cells.Item(33, 5).Value2 = '  +3.26%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value2 = '0.0680'
$ws.Cells.Item(34, 5).Value2 = '  +5.84%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value2 = '4.75'
$ws.Cells.Item(35, 5).Value2 = '  +3.96%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value2 = '4.85'
$ws.Cells.Item(36, 5).Value2 = '  -1.25%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value2 = '3.83'
$ws.Cells.Item(37, 5).Value2 = '  +9.15%  '

# Row 38
$ws.Cells.Item(38, 5).Value2 = '  +6.88%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value2 = '2.32'
$ws.Cells.Item(39, 5).Value2 = '  +0.06%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = '0.0257'
$ws.Cells.Item(40, 5).Value2 = '  +4.94%  '

# Row 41
$ws.Cells.Item(41, 5).Value2 = '  -0.14%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value2 = '17.49'
$ws.Cells.Item(42, 5).Value2 = '  +9.29%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value2 = '8.30'
$ws.Cells.Item(43, 5).Value2 = '  -1.96%  '

# Row 44
$ws.Cells.Item(44, 5).Value2 = '  +1.95%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value2 = '97.53'
$ws.Cells.Item(45, 5).Value2 = '  +1.20%  '

# Row 46: 'FTXToken' -> 'TerraClassic'
$ws.Cells.Item(46, 2).Value2 = 'TerraClassic'
$ws.Cells.Item(46, 3).Value2 = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value2 = '0.000213'
$ws.Cells.Item(46, 5).Value2 = '  -3.17%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value2 = '1.19'
$ws.Cells.Item(47, 5).Value2 = '  -0.09%  '

# Row 48: 'TerraClassic' -> 'FTXToken'
$ws.Cells.Item(48, 2).Value2 = 'FTXToken'
$ws.Cells.Item(48, 3).Value2 = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value2 = '4.37'
$ws.Cells.Item(48, 5).Value2 = '  +1.49%  '

# Row 49
$ws.Cells.Item(49, 4).Value2 = '1.445.67'
$ws.Cells.Item(49, 5).Value2 = '  -0.41%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value2 = '9.96'
$ws.Cells.Item(50, 5).Value2 = '  +18.57%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value2 = '2.30'
$ws.Cells.Item(51, 5).Value2 = '  +4.80%  '
